# Home page sliders testscript
#
# The locator sheet had a stray/duplicate row (row 255: "MainCategory_Xpath")
# sitting in the middle of the "Shopping Kart Page" section. Remove that
# entire row so the remaining rows shift up by one; this also drops the
# now-unused "MainCategory_Xpath" shared string from the string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("255").Delete()

# Leave the selection where it naturally lands after the delete (first
# empty row right after the last data row).
[void]$ws.Range("A280").Select()
